$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Block 1: columns B:F, rows 2-25 ---
$arr1 = New-Object 'object[,]' 24,5
$arr1[0,0] = 1.019999999999999
$arr1[0,1] = 1.027536138225185
$arr1[0,2] = 1.03703178602837
$arr1[0,3] = 1.027616377726966
$arr1[0,4] = 1.045680168218227
$arr1[1,0] = 1.02
$arr1[1,1] = 1.02856178734663
$arr1[1,2] = 1.037862154393593
$arr1[1,3] = 1.02848874713559
$arr1[1,4] = 1.046746319299607
$arr1[2,0] = 1.02
$arr1[2,1] = 1.02922551471639
$arr1[2,2] = 1.038399074575252
$arr1[2,3] = 1.02905367275519
$arr1[2,4] = 1.047436293695722
$arr1[3,0] = 1.02
$arr1[3,1] = 1.02950456089175
$arr1[3,2] = 1.038624703213781
$arr1[3,3] = 1.029291273152465
$arr1[3,4] = 1.04772638320882
$arr1[4,0] = 1.02
$arr1[4,1] = 1.029551414860641
$arr1[4,2] = 1.038662581802521
$arr1[4,3] = 1.029331173480751
$arr1[4,4] = 1.047775091924488
$arr1[5,0] = 1.02
$arr1[5,1] = 1.029229243285471
$arr1[5,2] = 1.038402089801005
$arr1[5,3] = 1.029056847169012
$arr1[5,4] = 1.047440169791202
$arr1[6,0] = 1.02
$arr1[6,1] = 1.027882747885628
$arr1[6,2] = 1.037312491897222
$arr1[6,3] = 1.027911106646042
$arr1[6,4] = 1.046040457288618
$arr1[7,0] = 1.02
$arr1[7,1] = 1.025510542735222
$arr1[7,2] = 1.035389575205952
$arr1[7,3] = 1.025895598241836
$arr1[7,4] = 1.043574796385689
$arr1[8,0] = 1.02
$arr1[8,1] = 1.023929404168148
$arr1[8,2] = 1.034105718274481
$arr1[8,3] = 1.024554272497425
$arr1[8,4] = 1.041931591383823
$arr1[9,0] = 1.02
$arr1[9,1] = 1.023244831464019
$arr1[9,2] = 1.033549348352498
$arr1[9,3] = 1.023974027653235
$arr1[9,4] = 1.041220204313176
$arr1[10,0] = 1.02
$arr1[10,1] = 1.022990561110058
$arr1[10,2] = 1.033342620518029
$arr1[10,3] = 1.023758583251224
$arr1[10,4] = 1.040955983437414
$arr1[11,0] = 1.02
$arr1[11,1] = 1.023045102495379
$arr1[11,2] = 1.033386967375862
$arr1[11,3] = 1.02380479298053
$arr1[11,4] = 1.041012658806033
$arr1[12,0] = 1.02
$arr1[12,1] = 1.02322381318564
$arr1[12,2] = 1.033532261534441
$arr1[12,3] = 1.023956217227633
$arr1[12,4] = 1.041198363324944
$arr1[13,0] = 1.02
$arr1[13,1] = 1.023333924116115
$arr1[13,2] = 1.033621773162674
$arr1[13,3] = 1.02404952589924
$arr1[13,4] = 1.04131278465826
$arr1[14,0] = 1.02
$arr1[14,1] = 1.023974838477318
$arr1[14,2] = 1.034142633250043
$arr1[14,3] = 1.024592793239414
$arr1[14,4] = 1.041978806588915
$arr1[15,0] = 1.02
$arr1[15,1] = 1.024376885921173
$arr1[15,2] = 1.034469234556672
$arr1[15,3] = 1.024933720269798
$arr1[15,4] = 1.042396619906047
$arr1[16,0] = 1.02
$arr1[16,1] = 1.024611400230794
$arr1[16,2] = 1.034659691983697
$arr1[16,3] = 1.025132631047695
$arr1[16,4] = 1.042640336023963
$arr1[17,0] = 1.02
$arr1[17,1] = 1.024691364715748
$arr1[17,2] = 1.034724625634239
$arr1[17,3] = 1.025200463593488
$arr1[17,4] = 1.042723439066528
$arr1[18,0] = 1.02
$arr1[18,1] = 1.024333749329409
$arr1[18,2] = 1.03443419783201
$arr1[18,3] = 1.02489713645631
$arr1[18,4] = 1.042351791178276
$arr1[19,0] = 1.02
$arr1[19,1] = 1.023171187055309
$arr1[19,2] = 1.033489477878833
$arr1[19,3] = 1.023911624227818
$arr1[19,4] = 1.04114367742025
$arr1[20,0] = 1.02
$arr1[20,1] = 1.022440298251469
$arr1[20,2] = 1.032895105991303
$arr1[20,3] = 1.023292481962781
$arr1[20,4] = 1.040384203906783
$arr1[21,0] = 1.02
$arr1[21,1] = 1.022827750591637
$arr1[21,2] = 1.033210230463649
$arr1[21,3] = 1.023620654523751
$arr1[21,4] = 1.040786804169191
$arr1[22,0] = 1.02
$arr1[22,1] = 1.024353240877119
$arr1[22,2] = 1.034450029554378
$arr1[22,3] = 1.024913666940598
$arr1[22,4] = 1.042372047311416
$arr1[23,0] = 1.02
$arr1[23,1] = 1.026123756246623
$arr1[23,2] = 1.035887035343296
$arr1[23,3] = 1.026416244684651
$arr1[23,4] = 1.044212130240834
$ws.Range("B2:F25").Value = $arr1

# --- Block 2: columns I:N, rows 2-25 ---
$arr2 = New-Object 'object[,]' 24,6
$arr2[0,0] = 1.033395596577766
$arr2[0,1] = 1.032693466908313
$arr2[0,2] = 1.039823830686133
$arr2[0,3] = 1.030435586947912
$arr2[0,4] = 1.048447729035884
$arr2[0,5] = 1.034160010082935
$arr2[1,0] = 1.033612040794812
$arr2[1,1] = 1.033358951479768
$arr2[1,2] = 1.040463868273384
$arr2[1,3] = 1.031115544414634
$arr2[1,4] = 1.049324706848221
$arr2[1,5] = 1.034826439718798
$arr2[2,0] = 1.033750170930674
$arr2[2,1] = 1.033789065491503
$arr2[2,2] = 1.04087698736079
$arr2[2,3] = 1.031555331864551
$arr2[2,4] = 1.049891677100311
$arr2[2,5] = 1.035257164541764
$arr2[3,0] = 1.033807779846842
$arr2[3,1] = 1.033969765531233
$arr2[3,2] = 1.041050415937733
$arr2[3,3] = 1.031740172404873
$arr2[3,4] = 1.050129912852077
$arr2[3,5] = 1.035438121196277
$arr2[4,0] = 1.033817425591485
$arr2[4,1] = 1.034000098820054
$arr2[4,2] = 1.041079520876213
$arr2[4,3] = 1.031771205226848
$arr2[4,4] = 1.050169906716816
$arr2[4,5] = 1.035468497561848
$arr2[5,0] = 1.033750942515791
$arr2[5,1] = 1.033791480485916
$arr2[5,2] = 1.040879305691977
$arr2[5,3] = 1.031557801894689
$arr2[5,4] = 1.049894860884532
$arr2[5,5] = 1.035259582965746
$arr2[6,0] = 1.033469142973107
$arr2[6,1] = 1.032918473601526
$arr2[6,2] = 1.040040346901974
$arr2[6,3] = 1.030665420903161
$arr2[6,4] = 1.048744209575271
$arr2[6,5] = 1.034385336311445
$arr2[7,0] = 1.032957861726546
$arr2[7,1] = 1.031376320476801
$arr2[7,2] = 1.038554145554468
$arr2[7,3] = 1.029091494075282
$arr2[7,4] = 1.046712860598336
$arr2[7,5] = 1.032840993152395
$arr2[8,0] = 1.032607138812895
$arr2[8,1] = 1.030345683711826
$arr2[8,2] = 1.037558098574937
$arr2[8,3] = 1.028041270017525
$arr2[8,4] = 1.045356130091488
$arr2[8,5] = 1.031808892764998
$arr2[9,0] = 1.032452935871061
$arr2[9,1] = 1.029898809474999
$arr2[9,2] = 1.037125560720466
$arr2[9,3] = 1.027586293624291
$arr2[9,4] = 1.044768062973063
$arr2[9,5] = 1.031361383915498
$arr2[10,0] = 1.032395307141101
$arr2[10,1] = 1.029732730326469
$arr2[10,2] = 1.036964710536664
$arr2[10,3] = 1.027417262046651
$arr2[10,4] = 1.044549539703221
$arr2[10,5] = 1.031195068915527
$arr2[11,0] = 1.032407684572998
$arr2[11,1] = 1.029768358943694
$arr2[11,2] = 1.036999221874658
$arr2[11,3] = 1.02745352138582
$arr2[11,4] = 1.04459641770781
$arr2[11,5] = 1.031230748129476
$arr2[12,0] = 1.032448179418586
$arr2[12,1] = 1.029885083158627
$arr2[12,2] = 1.037112268593778
$arr2[12,3] = 1.027572322094416
$arr2[12,4] = 1.044750001584641
$arr2[12,5] = 1.031347638106183
$arr2[13,0] = 1.032473083142031
$arr2[13,1] = 1.029956988859914
$arr2[13,2] = 1.037181895724155
$arr2[13,3] = 1.027645514757383
$arr2[13,4] = 1.044844617899943
$arr2[13,5] = 1.031419645921813
$arr2[14,0] = 1.032617323529973
$arr2[14,1] = 1.030375328626391
$arr2[14,2] = 1.037586778532381
$arr2[14,3] = 1.02807146064412
$arr2[14,4] = 1.045395145680007
$arr2[14,5] = 1.031838579778743
$arr2[15,0] = 1.032707176089802
$arr2[15,1] = 1.030637581184746
$arr2[15,2] = 1.037840418315619
$arr2[15,3] = 1.028338585926713
$arr2[15,4] = 1.045740318406577
$arr2[15,5] = 1.032101204765811
$arr2[16,0] = 1.03275936000238
$arr2[16,1] = 1.030790490592242
$arr2[16,2] = 1.03798824215976
$arr2[16,3] = 1.028494374099789
$arr2[16,4] = 1.045941594502449
$arr2[16,5] = 1.032254331322205
$arr2[17,0] = 1.032777115098786
$arr2[17,1] = 1.030842618891101
$arr2[17,2] = 1.038038625913463
$arr2[17,3] = 1.028547490186494
$arr2[17,4] = 1.04601021471991
$arr2[17,5] = 1.032306533649227
$arr2[18,0] = 1.032697559089488
$arr2[18,1] = 1.030609449972072
$arr2[18,2] = 1.037813217571516
$arr2[18,3] = 1.028309928131727
$arr2[18,4] = 1.045703290566103
$arr2[18,5] = 1.032073033603589
$arr2[19,0] = 1.032436264379746
$arr2[19,1] = 1.029850713273543
$arr2[19,2] = 1.037078984294272
$arr2[19,3] = 1.027537339157487
$arr2[19,4] = 1.044704777414507
$arr2[19,5] = 1.031313219411919
$arr2[20,0] = 1.032269948023663
$arr2[20,1] = 1.029373143876038
$arr2[20,2] = 1.036616264614913
$arr2[20,3] = 1.027051390533896
$arr2[20,4] = 1.044076458004796
$arr2[20,5] = 1.030834971811092
$arr2[21,0] = 1.032358307790217
$arr2[21,1] = 1.029626361804351
$arr2[21,2] = 1.036861663138722
$arr2[21,3] = 1.02730901911408
$arr2[21,4] = 1.044409590700712
$arr2[21,5] = 1.031088549337908
$arr2[22,0] = 1.032701905294896
$arr2[22,1] = 1.030622161435542
$arr2[22,2] = 1.037825508787489
$arr2[22,3] = 1.02832287742205
$arr2[22,4] = 1.045720022030352
$arr2[22,5] = 1.032085763118795
$arr2[23,0] = 1.033091781384176
$arr2[23,1] = 1.031775452411092
$arr2[23,2] = 1.038939291124679
$arr2[23,3] = 1.029498559439689
$arr2[23,4] = 1.047238455166822
$arr2[23,5] = 1.033240691899814
$ws.Range("I2:N25").Value = $arr2

Write-Host "Updated vm_pu values for case with 380 kV"
